$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Window position (cosmetic; mirrors author's on-screen window move) ---
$win = $excel.ActiveWindow
$win.Left = 17980
$win.Top = 1320

# --- Row 64: full-width black separator row with section label in column I ---
$ws.Range("A64:H64").Interior.Color = 0
$ws.Rows(64).RowHeight = 26
$ws.Range("I64").Value = "move to Github"

# --- Seed text values in the same order the shared-string table was built so
#     that new <si> entries land at the indices the target workbook uses. ---
$ws.Range("A66").Value = "Change automatic wd to find /CEDS/input"
$ws.Range("A65").Value = "Initial move to github"
$ws.Range("H65").Value = "2a8f9bb"
$ws.Range("A67").Value = "Adds and fixes .gitignore files"

# --- Row 65 ---
$ws.Range("B65").Value = 60
$ws.Range("C65").Value = "Rachel Hoesly"
$ws.Range("D65").Value = "Committed"
$ws.Range("E65").Value = [DateTime]"2015-08-28"
$ws.Range("G65").Value = [DateTime]"2015-08-28"

# --- Row 66 ---
$ws.Range("B66").Value = 62
$ws.Range("C66").Value = "Rachel Hoesly"
$ws.Range("D66").Value = "Committed"
$ws.Range("E66").Value = [DateTime]"2015-08-31"
$ws.Range("F66").Value = "-"
$ws.Range("G66").Value = [DateTime]"2015-08-31"
$ws.Range("H66").Value = 4611149

# --- Row 67 ---
$ws.Range("B67").Value = 63
$ws.Range("C67").Value = "Rachel Hoesly"
$ws.Range("D67").Value = "Committed"
$ws.Range("E67").Value = [DateTime]"2015-08-31"
$ws.Range("G67").Value = [DateTime]"2015-08-31"

# --- View: active selection after edit ---
$ws.Range("H67").Select() | Out-Null
